# Applies the #5 commit: adds proper column headers + company/insurance/debt
# metadata columns to the 保險 (Insurance) and 債務 (Debt) sheets, matching the
# 土地 (Land) sheet's layout (property_category/category/date/legislator_name/
# legislator_id/source_file/index).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value that must stay plain TEXT even when it looks like an
# ISO date (e.g. "2012-04-27"), which Excel would otherwise auto-convert to a
# date serial. We flip the cell to text format first, write the value, then
# copy the normal (non-text) cell format from a known-good neighbour back on
# top so the cell's number format matches the rest of the sheet.
function Set-TextValue($ws, $row, $col, $value, $formatSourceAddr) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ===========================================================================
# Sheet: 保險 (Insurance)
# ===========================================================================
$ws2 = $wb.Worksheets.Item("保險")

# -- Header row (row 1): B/C/D get corrected labels, E..K are brand new -----
$ws2.Cells.Item(1,2).Value = "company"
$ws2.Cells.Item(1,3).Value = "name"
$ws2.Cells.Item(1,4).Value = "owner"
$ws2.Cells.Item(1,5).Value = "property_category"
$ws2.Cells.Item(1,6).Value = "category"
$ws2.Cells.Item(1,7).Value = "date"
$ws2.Cells.Item(1,8).Value = "legislator_name"
$ws2.Cells.Item(1,9).Value = "legislator_id"
$ws2.Cells.Item(1,10).Value = "source_file"
$ws2.Cells.Item(1,11).Value = "index"

# header style (bold / bordered) for the newly added E1:K1 cells, copied
# from the existing D1 header cell
$ws2.Range("D1").Copy() | Out-Null
$ws2.Range("E1:K1").PasteSpecial(-4122) | Out-Null

# -- Data rows 2-5: B/C/D keep their values (only C3 / E4 / E5 actually
#    differ from what was already stored); E..K are new columns -----------
$insData = @(
    @{ Row=2; B="國泰人壽";     C="保本111終身";       Index=126 },
    @{ Row=3; B="國泰人壽";     C="美意年年終生壽險";   Index=127 },
    @{ Row=4; B="紐約國際人壽"; C="聚寶盆變額萬能壽險"; Index=128 },
    @{ Row=5; B="保誠人壽";     C="築夢人生";           Index=129 }
)

foreach ($r in $insData) {
    $row = $r.Row
    $ws2.Cells.Item($row,2).Value = $r.B
    $ws2.Cells.Item($row,3).Value = $r.C
    $ws2.Cells.Item($row,4).Value = "陳淑慧"
    $ws2.Cells.Item($row,5).Value = "insurance"
    $ws2.Cells.Item($row,6).Value = "normal"
    Set-TextValue $ws2 $row 7 "2012-04-27" "D$row"
    $ws2.Cells.Item($row,8).Value = "陳淑慧"
    $ws2.Cells.Item($row,9).Value = 1720
    $ws2.Cells.Item($row,10).Value = "tmpe56a1"
    $ws2.Cells.Item($row,11).Value = $r.Index
}

# ===========================================================================
# Sheet: 債務 (Debt)
# ===========================================================================
$ws3 = $wb.Worksheets.Item("債務")

# -- Header row (row 1): B/C get corrected labels, D..G are relabeled /
#    shifted, H..N are brand new ------------------------------------------
$ws3.Cells.Item(1,2).Value = "species"
$ws3.Cells.Item(1,3).Value = "debtor"
$ws3.Cells.Item(1,4).Value = "owner"
$ws3.Cells.Item(1,5).Value = "total"
$ws3.Cells.Item(1,6).Value = "register_date"
$ws3.Cells.Item(1,7).Value = "register_reason"
$ws3.Cells.Item(1,8).Value = "property_category"
$ws3.Cells.Item(1,9).Value = "category"
$ws3.Cells.Item(1,10).Value = "date"
$ws3.Cells.Item(1,11).Value = "legislator_name"
$ws3.Cells.Item(1,12).Value = "legislator_id"
$ws3.Cells.Item(1,13).Value = "source_file"
$ws3.Cells.Item(1,14).Value = "index"

$ws3.Range("C1").Copy() | Out-Null
$ws3.Range("H1:N1").PasteSpecial(-4122) | Out-Null

# -- Data rows 2-4 ------------------------------------------------------
$debtData = @(
    @{ Row=2; Debtor="元大商業銀行臺南市中西區民生路"; Owner="曾洋右等七人臺南市中西區中正路"; Total=7423744; RegDate="94年08月02日"; Index=151 },
    @{ Row=3; Debtor="元大商業銀行臺南市中西區民生路"; Owner="元大商業銀行臺南市中西區民生路"; Total=2367290; RegDate="96年10月15日"; Index=152 },
    @{ Row=4; Debtor="元大商業銀行臺南市中西區民生路"; Owner="安泰商業銀行臺南市中西區中山路"; Total=3067394; RegDate="94年08月02日"; Index=153 }
)

$ws3.Cells.Item(2,2).Value = "借款"
$ws3.Cells.Item(2,3).Value = "林南生"
$ws3.Cells.Item(2,4).Value = "曾洋右等七人臺南市中西區中正路"
$ws3.Cells.Item(2,5).Value = 7423744
$ws3.Cells.Item(2,6).Value = "94年08月02日"
$ws3.Cells.Item(2,7).Value = "借款"

$ws3.Cells.Item(3,2).Value = "借款"
$ws3.Cells.Item(3,3).Value = "林南生"
$ws3.Cells.Item(3,4).Value = "元大商業銀行臺南市中西區民生路"
$ws3.Cells.Item(3,5).Value = 2367290
$ws3.Cells.Item(3,6).Value = "96年10月15日"
$ws3.Cells.Item(3,7).Value = "借款"

$ws3.Cells.Item(4,2).Value = "借款"
$ws3.Cells.Item(4,3).Value = "林南生"
$ws3.Cells.Item(4,4).Value = "安泰商業銀行臺南市中西區中山路"
$ws3.Cells.Item(4,5).Value = 3067394
$ws3.Cells.Item(4,6).Value = "94年08月02日"
$ws3.Cells.Item(4,7).Value = "借款"

for ($row = 2; $row -le 4; $row++) {
    $ws3.Cells.Item($row,8).Value = "debt"
    $ws3.Cells.Item($row,9).Value = "normal"
    Set-TextValue $ws3 $row 10 "2012-04-27" "K$row"
    $ws3.Cells.Item($row,11).Value = "陳淑慧"
    $ws3.Cells.Item($row,12).Value = 1720
    $ws3.Cells.Item($row,13).Value = "tmpe56a1"
}
$ws3.Cells.Item(2,14).Value = 151
$ws3.Cells.Item(3,14).Value = 152
$ws3.Cells.Item(4,14).Value = 153

Write-Output "applied edits"
